# Gatchina city 2016 (added)
# Adds a new data row (row 60) for "Гатчина" (Gatchina), year 2016,
# right below the existing Gatchina rows (2019/2018/2017 in rows 57-59).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (styles) of the last existing row (59) down into the
# new row (60) so the new row inherits the same cell styles (centered
# alignment, yellow highlight on N/P/Q/R, etc.) without touching its values.
[void]$ws.Range("A59:U59").Copy()
[void]$ws.Range("A60:U60").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# name / year
$ws.Range("A60").Value = "Гатчина"
$ws.Range("B60").Value = 2016

# popsize / avgemployers / unemployed / avgsalary
$ws.Range("C60").Value = 95.178
$ws.Range("D60").Value = 21.987
$ws.Range("E60").Value = 230
$ws.Range("F60").Value = 38178.1

# livarea / beforeschool / docsperpop / bedsperpop / cliniccap (unknown)
$ws.Range("G60").Value = "???"
$ws.Range("H60").Value = "???"
$ws.Range("I60").Value = "???"
$ws.Range("J60").Value = "???"
$ws.Range("K60").Value = "???"

# invests
$ws.Range("L60").Value = 4169.8

# funds (unknown)
$ws.Range("M60").Value = "???"

# companies / factoriescap
$ws.Range("N60").Value = 4979
$ws.Range("O60").Value = 25080.5

# conscap / consnewareas / consnewapt (unknown)
$ws.Range("P60").Value = "???"
$ws.Range("Q60").Value = "???"
$ws.Range("R60").Value = "???"

# retailturnover / foodservturnover (computed from raw totals, like the rows above)
$ws.Range("S60").Formula = "=9719218.2/1000"
$ws.Range("T60").Formula = "= 123943/1000"

# saldo
$ws.Range("U60").Value = 3

# Match the author's final view state: scrolled near the new row, with R51 selected.
[void]$ws.Range("R51").Select()

Write-Output "Gatchina 2016 row added at row 60"
